{"js": "// Replace each three-digit-division-by-one-digit answer in the table\n// with its updated value, matching the diff exactly (old -> new).\nconst replacements = [\n  [\"691\u00f76=115, 1\", \"627\u00f74=156, 3\"],\n  [\"412\u00f78=51, 4\", \"732\u00f79=81, 3\"],\n  [\"368\u00f75=73, 3\", \"115\u00f73=38, 1\"],\n  [\"764\u00f76=127, 2\", \"646\u00f77=92, 2\"],\n  [\"426\u00f76=71, 0\", \"283\u00f78=35, 3\"],\n  [\"152\u00f72=76, 0\", \"326\u00f76=54, 2\"],\n  [\"296\u00f79=32, 8\", \"124\u00f74=31, 0\"],\n  [\"727\u00f74=181, 3\", \"416\u00f76=69, 2\"],\n  [\"711\u00f77=101, 4\", \"318\u00f77=45, 3\"],\n  [\"884\u00f79=98, 2\", \"890\u00f79=98, 8\"],\n  [\"774\u00f75=154, 4\", \"775\u00f75=155, 0\"],\n  [\"864\u00f75=172, 4\", \"249\u00f77=35, 4\"],\n  [\"740\u00f79=82, 2\", \"919\u00f73=306, 1\"],\n  [\"798\u00f72=399, 0\", \"886\u00f75=177, 1\"],\n  [\"697\u00f77=99, 4\", \"728\u00f77=104, 0\"],\n  [\"329\u00f72=164, 1\", \"379\u00f79=42, 1\"],\n  [\"958\u00f79=106, 4\", \"880\u00f73=293, 1\"],\n  [\"944\u00f75=188, 4\", \"982\u00f73=327, 1\"],\n  [\"267\u00f78=33, 3\", \"518\u00f77=74, 0\"],\n  [\"865\u00f79=96, 1\", \"303\u00f76=50, 3\"],\n  [\"174\u00f73=58, 0\", \"614\u00f78=76, 6\"],\n  [\"598\u00f73=199, 1\", \"822\u00f79=91, 3\"],\n  [\"561\u00f78=70, 1\", \"154\u00f78=19, 2\"],\n  [\"483\u00f77=69, 0\", \"597\u00f76=99, 3\"],\n  [\"619\u00f76=103, 1\", \"551\u00f75=110, 1\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const r of found.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# old answer -> new answer, as pairs (order matches the table, left-to-right, top-to-bottom)\n$replacements = @(\n    @(\"691\u00f76=115, 1\", \"627\u00f74=156, 3\"),\n    @(\"412\u00f78=51, 4\", \"732\u00f79=81, 3\"),\n    @(\"368\u00f75=73, 3\", \"115\u00f73=38, 1\"),\n    @(\"764\u00f76=127, 2\", \"646\u00f77=92, 2\"),\n    @(\"426\u00f76=71, 0\", \"283\u00f78=35, 3\"),\n    @(\"152\u00f72=76, 0\", \"326\u00f76=54, 2\"),\n    @(\"296\u00f79=32, 8\", \"124\u00f74=31, 0\"),\n    @(\"727\u00f74=181, 3\", \"416\u00f76=69, 2\"),\n    @(\"711\u00f77=101, 4\", \"318\u00f77=45, 3\"),\n    @(\"884\u00f79=98, 2\", \"890\u00f79=98, 8\"),\n    @(\"774\u00f75=154, 4\", \"775\u00f75=155, 0\"),\n    @(\"864\u00f75=172, 4\", \"249\u00f77=35, 4\"),\n    @(\"740\u00f79=82, 2\", \"919\u00f73=306, 1\"),\n    @(\"798\u00f72=399, 0\", \"886\u00f75=177, 1\"),\n    @(\"697\u00f77=99, 4\", \"728\u00f77=104, 0\"),\n    @(\"329\u00f72=164, 1\", \"379\u00f79=42, 1\"),\n    @(\"958\u00f79=106, 4\", \"880\u00f73=293, 1\"),\n    @(\"944\u00f75=188, 4\", \"982\u00f73=327, 1\"),\n    @(\"267\u00f78=33, 3\", \"518\u00f77=74, 0\"),\n    @(\"865\u00f79=96, 1\", \"303\u00f76=50, 3\"),\n    @(\"174\u00f73=58, 0\", \"614\u00f78=76, 6\"),\n    @(\"598\u00f73=199, 1\", \"822\u00f79=91, 3\"),\n    @(\"561\u00f78=70, 1\", \"154\u00f78=19, 2\"),\n    @(\"483\u00f77=69, 0\", \"597\u00f76=99, 3\"),\n    @(\"619\u00f76=103, 1\", \"551\u00f75=110, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n"}
